# Second commit - updated Dashboard page, DashboardPageTest
# Adds a new "DashboardPageTabs" worksheet listing the tabs shown on the
# Salesforce dashboard page, with a green header cell.

$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item(1)

# Add the new worksheet after the existing (last) sheet and name it.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "DashboardPageTabs"

# Base list of dashboard page tabs, typed first.
$base = @(
    "dashboardPageTabs",
    "Home",
    "Accounts",
    "Contacts",
    "Leads",
    "Tasks",
    "Calendar",
    "Dashboards",
    "Reports",
    "Groups",
    "Forecasts",
    "Files",
    "Quotes",
    "Chatter"
)

for ($i = 0; $i -lt $base.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $base[$i]
}

# A few tabs were added afterwards, inserted in their respective spot
# in the list (pushing the following rows down).
$ws.Rows.Item(6).Insert(-4121) | Out-Null
$ws.Cells.Item(6, 1).Value = "Opportunities"

$ws.Rows.Item(14).Insert(-4121) | Out-Null
$ws.Cells.Item(14, 1).Value = "List Emails"

$ws.Rows.Item(10).Insert(-4121) | Out-Null
$ws.Cells.Item(10, 1).Value = "Notes"

# Reuse the existing bordered/no-fill style (same one already used on
# LoginTestData!B2) for all the plain data rows A2:A17.
$loginSheet.Range("B2").Copy() | Out-Null
$ws.Range("A2:A17").PasteSpecial(-4122) | Out-Null

# Reuse the existing bold+border header style for A1, then recolor its
# fill to green, creating the new header style used for the tab title.
$loginSheet.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Interior.Color = 5287936

$ws.Application.CutCopyMode = $false

# Column A width, fit to the longest tab name.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Match the selection/active cell recorded for the new sheet and make it
# the active tab of the workbook.
$ws.Range("A10").Select() | Out-Null
